# Natmi following Dr Hou advice
# Rewrite the LR-pair data rows (2-4 updated, 5-7 added) with the new
# sending/target cluster assignments and recomputed edge statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
  @("ECs", "Cdh1", "Igf1r", "ECs", 1, 0.3333333333333333, 0.098866, 0.296598, 0.3026185969870575, 0.3026185969870575, 3, 1, 29.57110033333333, 88.713301, 0.5311523066901919, 0.5311523066901919, 2.923576405555333, 26.312187649998, 0.1607365658370251, 0.1607365658370251),
  @("ECs", "Cdh1", "Igf1r", "FAPs", 1, 0.3333333333333333, 0.098866, 0.296598, 0.3026185969870575, 0.3026185969870575, 3, 1, 12.180337, 36.541011, 0.2187816490048282, 0.2187816490048282, 1.204221197842, 10.837990780578, 0.06620739566835598, 0.06620739566835597),
  @("ECs", "Cdh1", "Igf1r", "sCs", 1, 0.3333333333333333, 0.098866, 0.296598, 0.3026185969870575, 0.3026185969870575, 3, 1, 13.92204833333333, 41.76614499999999, 0.2500660443049799, 0.2500660443049799, 1.376417230523333, 12.38775507471, 0.07567463548167638, 0.07567463548167637),
  @("sCs", "Cdh1", "Igf1r", "ECs", 3, 1, 0.2278356666666667, 0.6835070000000001, 0.6973814030129426, 0.6973814030129426, 3, 1, 29.57110033333333, 88.713301, 0.5311523066901919, 0.5311523066901919, 6.73735135851189, 60.63616222660701, 0.3704157408531668, 0.3704157408531668),
  @("sCs", "Cdh1", "Igf1r", "FAPs", 3, 1, 0.2278356666666667, 0.6835070000000001, 0.6973814030129426, 0.6973814030129426, 3, 1, 12.180337, 36.541011, 0.2187816490048282, 0.2187816490048282, 2.775115200619667, 24.97603680557701, 0.1525742533364723, 0.1525742533364723),
  @("sCs", "Cdh1", "Igf1r", "sCs", 3, 1, 0.2278356666666667, 0.6835070000000001, 0.6973814030129426, 0.6973814030129426, 3, 1, 13.92204833333333, 41.76614499999999, 0.2500660443049799, 0.2500660443049799, 3.171939163390555, 28.547452470515, 0.1743914088233036, 0.1743914088233036)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $rows[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
